# Add the new "2022-Q4" detail sheet right after "总计" (i.e. before "2022-Q3"),
# and add the corresponding summary row to "总计".
#
# All existing "20XX-Qn" tabs shift right by one position; their own content
# is untouched (the engine assigns worksheets/sheetN.xml + rIds by position
# automatically), so we only need to (a) insert+populate the new sheet and
# (b) edit the "总计" sheet in place.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q4, push the rest
#    down, and bump their index column (A) by one.
#
# Avoid Rows.Insert here - it inherits formatting from the row above
# (the bold header), which would mint an unwanted blended style. Instead
# shift the existing rows down manually (bottom-up) with plain value
# writes, which leaves every cell's style exactly as it already was.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $summary.Range("A$dst").Value = $summary.Range("A$src").Value
    $summary.Range("B$dst").Value = $summary.Range("B$src").Value
    $summary.Range("C$dst").Value = $summary.Range("C$src").Value
    $summary.Range("D$dst").Value = $summary.Range("D$src").Value
}

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 1.18

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" detail sheet before the current "2022-Q3"
#    (i.e. right after "总计", in slot 2).
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$new.Name = "2022-Q4"

# Match the page margins used by every other detail sheet.
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36

# A cell elsewhere in the workbook that already carries the bold/bordered/
# centered style (cellXfs index 2) used for headers + the index column -
# copy-pasting its format lets us reuse that existing style instead of
# minting a new one.
$styledCell = $wb.Worksheets.Item(3).Range("B1")
$styledCell.Copy() | Out-Null

# ---- header row ----
$new.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# ---- index column (A2:A10) ----
$styledCell.Copy() | Out-Null
$new.Range("A2:A10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
for ($i = 0; $i -le 8; $i++) {
    $new.Range("A$($i + 2)").Value = $i
}

# ---- data rows ----
$rows = @(
    @("009190", "景顺长城核心优选一年持有期混合", "10.80", "90.64", "3.24", "0.3499", 10),
    @("005888", "华夏新兴消费混合A",               "7.99",  "88.30", "3.23", "0.2581", 4),
    @("012421", "华夏优加生活混合A",               "8.01",  "88.08", "3.14", "0.2515", 4),
    @("005889", "华夏新兴消费混合C",               "4.95",  "88.30", "3.23", "0.1599", 4),
    @("001703", "银华沪港深增长股票A",             "1.71",  "93.53", "4.97", "0.0850", 5),
    @("015797", "万家新能源主题混合C",             "0.69",  "93.18", "5.77", "0.0398", 2),
    @("015796", "万家新能源主题混合A",             "0.31",  "93.18", "5.77", "0.0179", 2),
    @("014364", "银华沪港深增长股票C",             "0.25",  "93.53", "4.97", "0.0124", 5),
    @("012422", "华夏优加生活混合C",               "0.32",  "88.08", "3.14", "0.0100", 4)
)

# Columns B, D, E, F, G hold numeric-looking text (fund codes / percentages /
# values) that must stay text (leading zeros, trailing zeros, 4dp precision).
# Force-format as text first so the values aren't silently parsed as numbers,
# then clear the formatting again so the cells end up with the default style
# (matching the source file, which never touches styles.xml for these cells).
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $colRange = $new.Range("$col" + "2:" + "$col" + "10")
    $colRange.NumberFormat = "@"
}

for ($i = 0; $i -le 8; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    $new.Range("B$r").Value = $vals[0]
    $new.Range("C$r").Value = $vals[1]
    $new.Range("D$r").Value = $vals[2]
    $new.Range("E$r").Value = $vals[3]
    $new.Range("F$r").Value = $vals[4]
    $new.Range("G$r").Value = $vals[5]
    $new.Range("H$r").Value = $vals[6]
}

foreach ($col in $textCols) {
    $colRange = $new.Range("$col" + "2:" + "$col" + "10")
    $colRange.ClearFormats()
}
